$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Custos"
$ws.Range("B7").Value = "R$ 154.000"

$ws.Range("B8").Select()
